# "edit inactive ads complete"
# Score the "Edit Inactive Ads" self-evaluation rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the completed scores (B26 = "Edit Inactive Ads", B27 = "Change and Edit Images").
$ws.Range("C26").Value = 10
$ws.Range("C27").Value = 5

# Recalculate so the dependent SUM formula (C51) picks up the new total.
$excel.Calculate()

# Move the active selection to where the user ended up editing.
$ws.Range("E26").Select()
